# Update the lattice-multiplication exercise table: each of the 15 cells
# (5 rows x 3 columns) keeps its layout (problem header, factors, rule,
# two partial-product rows separated by line breaks) but gets new numbers.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$vt = [char]11   # vertical-tab == the in-cell line break Word exposes for <w:br/>

function Set-CellLines($table, $row, $col, $lines) {
    $table.Cell($row, $col).Range.Text = [string]::Join($vt, $lines)
}

Set-CellLines $t 1 1 @("96 x 53", "  5    3", "  ----", "9|    |", "6|    |")
Set-CellLines $t 1 2 @("99 x 82", "  8    2", "  ----", "9|    |", "9|    |")
Set-CellLines $t 1 3 @("17 x 76", "  7    6", "  ----", "1|    |", "7|    |")

Set-CellLines $t 2 1 @("26 x 10", "  1    0", "  ----", "2|    |", "6|    |")
Set-CellLines $t 2 2 @("77 x 66", "  6    6", "  ----", "7|    |", "7|    |")
Set-CellLines $t 2 3 @("28 x 13", "  1    3", "  ----", "2|    |", "8|    |")

Set-CellLines $t 3 1 @("63 x 48", "  4    8", "  ----", "6|    |", "3|    |")
Set-CellLines $t 3 2 @("80 x 37", "  3    7", "  ----", "8|    |", "0|    |")
Set-CellLines $t 3 3 @("20 x 42", "  4    2", "  ----", "2|    |", "0|    |")

Set-CellLines $t 4 1 @("51 x 83", "  8    3", "  ----", "5|    |", "1|    |")
Set-CellLines $t 4 2 @("30 x 63", "  6    3", "  ----", "3|    |", "0|    |")
Set-CellLines $t 4 3 @("68 x 62", "  6    2", "  ----", "6|    |", "8|    |")

Set-CellLines $t 5 1 @("82 x 88", "  8    8", "  ----", "8|    |", "2|    |")
Set-CellLines $t 5 2 @("36 x 25", "  2    5", "  ----", "3|    |", "6|    |")
Set-CellLines $t 5 3 @("55 x 27", "  2    7", "  ----", "5|    |", "5|    |")
